$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 header changes
$ws.Range("K3").Value = "P(x0|Z)"
$ws.Range("H3").ClearContents()
$ws.Range("L3").ClearContents()

# Row 2 label changes
$ws.Range("B2").Value = "Medicine"
$ws.Range("F2").Value = "Placebo"

# Row 4 values
$ws.Range("G4").Value = 0.49
$ws.Range("K4").Value = 0.61
$ws.Range("H4").ClearContents()
$ws.Range("L4").ClearContents()

# Row 5 values
$ws.Range("G5").Value = 0.08
$ws.Range("K5").Value = 0.04
$ws.Range("H5").ClearContents()
$ws.Range("L5").ClearContents()

# Update selection to L4 as shown in diff
$ws.Range("L4").Select()
